$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.542.91'
$ws.Range("E2").Value = '  +5.52%  '
$ws.Range("D3").Value = '1.723.52'
$ws.Range("E3").Value = '  +4.21%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '225.50'
$ws.Range("E5").Value = '  +3.30%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5365'
$ws.Range("E6").Value = '  +3.14%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2671'
$ws.Range("E8").Value = '  +0.96%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06600'
$ws.Range("E9").Value = '  +4.15%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.73'
$ws.Range("E10").Value = '  +6.40%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07726'
$ws.Range("E11").Value = '  +0.47%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.611'
$ws.Range("E12").Value = '  +0.15%  '
$ws.Range("D13").Value = '1.721.85'
$ws.Range("E13").Value = '  +5.70%  '
$ws.Range("D14").Value = '1.961.11'
$ws.Range("E14").Value = '  +4.31%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5851'
$ws.Range("E15").Value = '  +4.57%  '
$ws.Range("D16").Value = '0.0₅8309'
$ws.Range("E16").Value = '  +1.94%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.98'
$ws.Range("E17").Value = '  +3.97%  '
$ws.Range("D18").Value = '27.552.92'
$ws.Range("E18").Value = '  +5.52%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '220.41'
$ws.Range("E19").Value = '  +15.22%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.003'
$ws.Range("E20").Value = '  +0.00%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.723'
$ws.Range("E21").Value = '  +1.96%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.65'
$ws.Range("E22").Value = '  +1.55%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.090'
$ws.Range("E23").Value = '  +2.79%  '
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("E25").Value = '  +3.11%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.742'
$ws.Range("E26").Value = '  +15.91%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1236'
$ws.Range("E27").Value = '  +3.89%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.418'
$ws.Range("E28").Value = '  +2.71%  '
$ws.Range("E29").Value = '  +4.82%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05600'
$ws.Range("E30").Value = '  +2.03%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.303'
$ws.Range("E31").Value = '  +2.69%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.556'
$ws.Range("E32").Value = '  +3.27%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.453'
$ws.Range("E33").Value = '  +2.88%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.660'
$ws.Range("E34").Value = '  +6.46%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.835'
$ws.Range("E35").Value = '  +1.79%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9608'
$ws.Range("E36").Value = '  +1.36%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.428'
$ws.Range("E37").Value = '  +0.23%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5962'
$ws.Range("E38").Value = '  +5.67%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01652'
$ws.Range("E39").Value = '  +4.57%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.929'
$ws.Range("E40").Value = '  +1.42%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8542'
$ws.Range("E41").Value = '  +3.12%  '
$ws.Range("D42").Value = '1.056.02'
$ws.Range("E42").Value = '  +2.77%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.003'
$ws.Range("E43").Value = '  +0.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.66'
$ws.Range("E44").Value = '  +0.54%  '
$ws.Range("D45").Value = '1.867.44'
$ws.Range("E45").Value = '  +4.03%  '
$ws.Range("E46").Value = '  +5.96%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '59.05'
$ws.Range("E47").Value = '  +2.53%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.215'
$ws.Range("E48").Value = '  +3.04%  '
$ws.Range("E49").Value = '  +2.28%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05253'
$ws.Range("E51").Value = '  +1.57%  '
